# Update gh-pages output numbers (想去人数 / 最低票价) for 苏州-漫展信息.xlsx
# Applies to both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- "展览" sheet (sheet1) ---
$ws1.Range("F2").Value = 3159
$ws1.Range("F5").Value = 99
$ws1.Range("F6").Value = 53
$ws1.Range("F10").Value = 16001
$ws1.Range("F11").Value = 257
$ws1.Range("F14").Value = 6257
$ws1.Range("F15").Value = 632
$ws1.Range("F16").Value = 118
$ws1.Range("F20").Value = 1266
$ws1.Range("F29").Value = 5023
$ws1.Range("G29").Value = 68
$ws1.Range("F31").Value = 11196
$ws1.Range("F34").Value = 134
$ws1.Range("F35").Value = 188
$ws1.Range("F36").Value = 3820

# --- "全部类型" sheet (sheet4) ---
$ws4.Range("F2").Value = 3159
$ws4.Range("F5").Value = 99
$ws4.Range("F6").Value = 53
$ws4.Range("F10").Value = 16001
$ws4.Range("F11").Value = 257
$ws4.Range("F14").Value = 6257
$ws4.Range("F15").Value = 632
$ws4.Range("F16").Value = 118
$ws4.Range("F20").Value = 1266
$ws4.Range("F29").Value = 5023
$ws4.Range("G29").Value = 68
$ws4.Range("F32").Value = 11196
$ws4.Range("F35").Value = 134
$ws4.Range("F36").Value = 188
$ws4.Range("F37").Value = 3820
